# Apply the "Atualizacao de bases das ligas" edit:
#  - Several same-date fixtures had their data rows cyclically re-shuffled
#    (row numbers / A-column sequence IDs stay fixed; B:AC contents rotate).
#  - Three brand-new match rows (201-203) were appended at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-RowValues {
    param($ws, $rangeAddr, $values)
    $n = $values.Count
    $arr = New-Object "object[,]" 1,$n
    for ($i = 0; $i -lt $n; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($rangeAddr).Value = $arr
}

# --- Update rows whose data (columns B:AC) was cyclically shifted within same-date groups ---
Set-RowValues $ws "B26:AC26" @(6156883, "Chile Primera B", "Chile Primera B", 45121.85416666666, "Universidad de Concepcion", "Deportes Recoleta", 1, 1, "D", 2.2, 3, 3.1, 1.85, 3.2, 3.75, -0.5, 1.925, 1.875, 2.5, 1.825, 1.975, -1, 2.2, -1, -1, 0.875, -1, 0.9750000000000001)
Set-RowValues $ws "B27:AC27" @(6156881, "Chile Primera B", "Chile Primera B", 45121.85416666666, "San Luis Quillota", "Deportes Iquique", 0, 1, "A", 2.1, 3.2, 3.2, 2.05, 3.2, 3.3, -0.25, 1.8, 2, 2.5, 1.95, 1.85, -1, -1, 2.3, -1, 1, -1, 0.8500000000000001)
Set-RowValues $ws "B35:AC35" @(6156888, "Chile Primera B", "Chile Primera B", 45129.5625, "Deportes Temuco", "San Luis Quillota", 1, 2, "A", 2.2, 2.875, 3.2, 2.15, 2.9, 3.25, -0.25, 1.9, 1.9, 2.25, 1.95, 1.85, -1, -1, 2.25, -1, 0.8999999999999999, 0.95, -1)
Set-RowValues $ws "B36:AC36" @(6156886, "Chile Primera B", "Chile Primera B", 45129.5625, "Union San Felipe", "San Marcos De Arica", 3, 0, "H", 2.15, 3, 3.2, 2.1, 3, 3.25, -0.25, 1.875, 1.925, 2.5, 1.875, 1.925, 1.1, -1, -1, 0.875, -1, 0.875, -1)
Set-RowValues $ws "B38:AC38" @(6158713, "Chile Primera B", "Chile Primera B", 45130.5625, "Santiago Morning", "CD Antofagasta", 1, 2, "A", 3, 3.25, 2.1, 2.625, 3.2, 2.3, 0, 2.025, 1.775, 2.25, 1.75, 1.95, -1, -1, 1.3, -1, 0.7749999999999999, 0.75, -1)
Set-RowValues $ws "B39:AC39" @(6155944, "Chile Primera B", "Chile Primera B", 45130.5625, "Cobreloa", "Puerto Montt", 1, 0, "H", 1.4, 3.75, 7.5, 1.45, 3.6, 6.5, -1.25, 2.05, 1.75, 2.25, 1.9, 1.9, 0.45, -1, -1, -0.5, 0.375, -1, 0.8999999999999999)
Set-RowValues $ws "B43:AC43" @(6393620, "Chile Primera B", "Chile Primera B", 45136.5625, "Union San Felipe", "Santiago Morning", 2, 1, "H", 1.95, 3, 3.75, 1.65, 3.25, 5, -0.75, 1.9, 1.9, 2.5, 1.9, 1.9, 0.6499999999999999, -1, -1, 0.45, -0.5, 0.8999999999999999, -1)
Set-RowValues $ws "B44:AC44" @(6155945, "Chile Primera B", "Chile Primera B", 45136.5625, "Puerto Montt", "Deportes Iquique", 2, 2, "D", 2.625, 3.2, 2.375, 2.6, 3.25, 2.375, 0, 2, 1.8, 2.25, 1.75, 1.95, -1, 2.25, -1, 0, -0, 0.75, -1)
Set-RowValues $ws "B69:AC69" @(6156911, "Chile Primera B", "Chile Primera B", 45158.5625, "Barnechea", "Universidad de Concepcion", 2, 0, "H", 2.1, 3.2, 3.1, 2.15, 3.2, 3, -0.25, 1.975, 1.825, 2.5, 1.875, 1.925, 1.15, -1, -1, 0.9750000000000001, -1, -1, 0.925)
Set-RowValues $ws "B70:AC70" @(6156908, "Chile Primera B", "Chile Primera B", 45158.5625, "Deportes Temuco", "Union San Felipe", 2, 1, "H", 2.3, 3.2, 2.7, 2, 3.25, 3.2, -0.25, 1.8, 2, 2.5, 1.9, 1.9, 1, -1, -1, 0.8, -1, 0.8999999999999999, -1)
Set-RowValues $ws "B101:AC101" @(6155958, "Chile Primera B", "Chile Primera B", 45193.52083333334, "Puerto Montt", "San Marcos De Arica", 1, 1, "D", 3, 3.2, 2.2, 2.05, 3.3, 3.25, -0.25, 1.825, 1.975, 2.5, 1.95, 1.75, -1, 2.3, -1, -0.5, 0.4875, -1, 0.75)
Set-RowValues $ws "B102:AC102" @(6156183, "Chile Primera B", "Chile Primera B", 45193.52083333334, "Santiago Wanderers", "CD Antofagasta", 3, 0, "H", 2.2, 3.2, 3, 2.2, 3.2, 3, -0.25, 1.95, 1.85, 2.25, 1.8, 2, 1.2, -1, -1, 0.95, -1, 0.8, -1)
Set-RowValues $ws "B103:AC103" @(6156929, "Chile Primera B", "Chile Primera B", 45193.52083333334, "Cobreloa", "Union San Felipe", 1, 0, "H", 2.3, 2.9, 3, 1.8, 3.3, 4, -0.5, 1.875, 1.925, 2.5, 1.9, 1.9, 0.8, -1, -1, 0.875, -1, -1, 0.8999999999999999)
Set-RowValues $ws "B114:AC114" @(6156940, "Chile Primera B", "Chile Primera B", 45207.72916666666, "Universidad de Concepcion", "Union San Felipe", 0, 0, "D", 2.3, 3.1, 2.9, 2.2, 3.2, 3.1, -0.25, 1.9, 1.9, 2.5, 1.825, 1.975, -1, 2.2, -1, -0.5, 0.45, -1, 0.9750000000000001)
Set-RowValues $ws "B115:AC115" @(6155961, "Chile Primera B", "Chile Primera B", 45207.72916666666, "La Serena", "Santiago Morning", 0, 1, "A", 1.8, 3.3, 4, 2.05, 3.2, 3.3, -0.5, 2.05, 1.75, 2.5, 1.925, 1.875, -1, -1, 2.3, -1, 0.75, -1, 0.875)
Set-RowValues $ws "B118:AC118" @(6156939, "Chile Primera B", "Chile Primera B", 45208.72916666666, "Club Deportes Santa Cruz", "Deportes Iquique", 1, 2, "A", 2.375, 3.1, 2.75, 2.45, 3.2, 2.625, 0, 1.775, 2.025, 2.5, 1.85, 1.95, -1, -1, 1.625, -1, 1.025, 0.8500000000000001, -1)
Set-RowValues $ws "B120:AC120" @(6156943, "Chile Primera B", "Chile Primera B", 45208.72916666666, "Santiago Wanderers", "Deportes Temuco", 1, 0, "H", 2.1, 3.2, 3.1, 1.909, 3.4, 3.4, -0.5, 1.95, 1.85, 2.25, 1.825, 1.975, 0.909, -1, -1, 0.95, -1, -1, 0.9750000000000001)
Set-RowValues $ws "B122:AC122" @(7327856, "Chile Primera B", "Chile Primera B", 45213.70833333334, "Union San Felipe", "Puerto Montt", 0, 1, "A", 1.727, 3.5, 4, 1.8, 3.4, 3.75, -0.5, 1.85, 1.95, 2.25, 1.75, 1.95, -1, -1, 2.75, -1, 0.95, -1, 0.95)
Set-RowValues $ws "B123:AC123" @(7327855, "Chile Primera B", "Chile Primera B", 45213.70833333334, "Santiago Morning", "Deportes Recoleta", 0, 1, "A", 2.1, 3.4, 3, 2.625, 3.3, 2.4, 0, 2, 1.8, 2.5, 1.85, 1.95, -1, -1, 1.4, -1, 0.8, -1, 0.95)
Set-RowValues $ws "B124:AC124" @(7327838, "Chile Primera B", "Chile Primera B", 45213.70833333334, "Barnechea", "San Marcos De Arica", 3, 3, "D", 2, 3.3, 3.2, 2.1, 3.3, 3, -0.25, 1.9, 1.9, 3, 2, 1.8, -1, 2.3, -1, -0.5, 0.45, 1, -1)
Set-RowValues $ws "B125:AC125" @(7327842, "Chile Primera B", "Chile Primera B", 45214.52083333334, "Deportes Temuco", "Club Deportes Santa Cruz", 2, 1, "H", 1.533, 4, 5, 1.7, 3.75, 4, -0.75, 1.9, 1.9, 2.25, 1.825, 1.975, 0.7, -1, -1, 0.45, -0.5, 0.825, -1)
Set-RowValues $ws "B127:AC127" @(7327843, "Chile Primera B", "Chile Primera B", 45214.52083333334, "San Luis Quillota", "Universidad de Concepcion", 2, 3, "A", 1.727, 3.6, 4, 1.727, 3.4, 4.2, -0.75, 1.975, 1.825, 2.5, 1.925, 1.875, -1, -1, 3.2, -1, 0.825, 0.925, -1)
Set-RowValues $ws "B128:AC128" @(7327840, "Chile Primera B", "Chile Primera B", 45214.52083333334, "CD Antofagasta", "La Serena", 0, 1, "A", 1.8, 3.3, 4, 2, 2.9, 3.75, -0.25, 1.725, 1.975, 2.5, 1.85, 1.95, -1, -1, 2.75, -1, 0.9750000000000001, -1, 0.95)
Set-RowValues $ws "B129:AC129" @(7327841, "Chile Primera B", "Chile Primera B", 45214.52083333334, "Deportes Iquique", "Santiago Wanderers", 3, 2, "H", 2.1, 3.2, 3.2, 2.05, 3.2, 3.25, -0.25, 1.825, 1.975, 2.5, 1.925, 1.875, 1.05, -1, -1, 0.825, -1, 0.925, -1)
Set-RowValues $ws "B136:AC136" @(7487575, "Chile Primera B", "Chile Primera B", 45249.75, "CD Antofagasta", "Deportes Iquique", 2, 2, "D", 2.1, 3.5, 3, 1.95, 3.6, 3.4, -0.5, 1.975, 1.825, 2.75, 1.975, 1.825, -1, 2.6, -1, -1, 0.825, 0.9750000000000001, -1)
Set-RowValues $ws "B137:AC137" @(7487574, "Chile Primera B", "Chile Primera B", 45249.75, "Deportes Temuco", "Santiago Wanderers", 2, 3, "A", 2, 3.2, 3.4, 1.95, 3.2, 3.5, -0.5, 2, 1.8, 2.5, 1.975, 1.825, -1, -1, 2.5, -1, 0.8, 0.9750000000000001, -1)
Set-RowValues $ws "B138:AC138" @(7503210, "Chile Primera B", "Chile Primera B", 45256.75, "Santiago Wanderers", "Deportes Temuco", 1, 0, "H", 2.2, 3.2, 3, 2.15, 3.2, 3.1, -0.25, 1.925, 1.875, 2.5, 1.975, 1.825, 1.15, -1, -1, 0.925, -1, -1, 0.825)
Set-RowValues $ws "B139:AC139" @(7503209, "Chile Primera B", "Chile Primera B", 45256.75, "Deportes Iquique", "CD Antofagasta", 3, 0, "H", 2.2, 3.6, 2.75, 2.15, 3.6, 2.875, -0.25, 1.975, 1.825, 2.75, 1.85, 1.95, 1.15, -1, -1, 0.9750000000000001, -1, 0.425, -0.5)
Set-RowValues $ws "B156:AC156" @(7794219, "Chile Primera B", "Chile Primera B", 45359.75, "Deportes Limache", "Rangers de Talca", 0, 1, "A", 2.2, 3.3, 2.9, 2.15, 3.3, 3, -0.25, 1.95, 1.85, 2.5, 1.975, 1.825, -1, -1, 2, -1, 0.8500000000000001, -1, 0.825)
Set-RowValues $ws "B157:AC157" @(7793612, "Chile Primera B", "Chile Primera B", 45359.75, "Deportes Recoleta", "Barnechea", 1, 2, "A", 1.909, 3.3, 3.75, 2.75, 3.4, 2.3, 0.25, 1.75, 2.05, 2.75, 1.975, 1.825, -1, -1, 1.3, -1, 1.05, 0.4875, -0.5)
Set-RowValues $ws "B158:AC158" @(7793492, "Chile Primera B", "Chile Primera B", 45359.85416666666, "Curico Unido", "Club Deportes Santa Cruz", 1, 2, "A", 1.8, 3.4, 4, 2.3, 3.1, 2.9, -0.25, 2.025, 1.775, 2.5, 1.975, 1.825, -1, -1, 1.9, -1, 0.7749999999999999, 0.9750000000000001, -1)
Set-RowValues $ws "B159:AC159" @(7793493, "Chile Primera B", "Chile Primera B", 45359.85416666666, "San Marcos De Arica", "Santiago Morning", 2, 0, "H", 2, 3.3, 3.3, 1.85, 3.5, 3.6, -0.5, 1.925, 1.875, 2.75, 2, 1.8, 0.8500000000000001, -1, -1, 0.925, -1, -1, 0.8)
Set-RowValues $ws "B177:AC177" @(7792884, "Chile Primera B", "Chile Primera B", 45375.75, "Santiago Wanderers", "Deportes Temuco", 2, 1, "H", 1.727, 3.6, 4, 1.727, 3.3, 4.2, -0.5, 1.8, 2, 2.5, 1.95, 1.75, 0.7270000000000001, -1, -1, 0.8, -1, 0.95, -1)
Set-RowValues $ws "B178:AC178" @(7793501, "Chile Primera B", "Chile Primera B", 45375.75, "Magallanes", "Santiago Morning", 2, 1, "H", 2.3, 3.4, 2.625, 1.909, 3.6, 3.3, -0.5, 1.925, 1.875, 2.5, 1.9, 1.9, 0.909, -1, -1, 0.925, -1, 0.8999999999999999, -1)
Set-RowValues $ws "B191:AC191" @(7793508, "Chile Primera B", "Chile Primera B", 45389.5625, "Deportes Recoleta", "San Marcos De Arica", 1, 2, "A", 2.15, 3.3, 2.875, 2.05, 3.4, 3, -0.25, 1.85, 1.95, 2.5, 1.9, 1.9, -1, -1, 2, -1, 0.95, 0.8999999999999999, -1)
Set-RowValues $ws "B192:AC192" @(7792888, "Chile Primera B", "Chile Primera B", 45389.5625, "CD Antofagasta", "Union San Felipe", 4, 1, "H", 1.75, 3.4, 4, 1.75, 3.4, 4, -0.5, 1.8, 2, 2.5, 1.875, 1.925, 0.75, -1, -1, 0.8, -1, 0.875, -1)

# --- Append new rows 201-203 ---
$ws.Cells.Item(201, 1).Value = 199
$ws.Cells.Item(201, 1).Style = $ws.Cells.Item(200, 1).Style
Set-RowValues $ws "B201:AC201" @(7793513, "Chile Primera B", "Chile Primera B", 45396.77083333334, "San Luis Quillota", "La Serena", 0, 1, "A", 2.75, 3.4, 2.2, 1.8, 3.6, 3.5, -0.5, 1.85, 1.95, 2.25, 1.85, 1.95, -1, -1, 2.5, -1, 0.95, -1, 0.95)
$ws.Cells.Item(201, 5).Style = $ws.Cells.Item(200, 5).Style

$ws.Cells.Item(202, 1).Value = 200
$ws.Cells.Item(202, 1).Style = $ws.Cells.Item(200, 1).Style
Set-RowValues $ws "B202:AC202" @(7792892, "Chile Primera B", "Chile Primera B", 45396.77083333334, "Union San Felipe", "Deportes Temuco", 2, 3, "A", 2.375, 3.3, 2.625, 2.1, 3.2, 3.1, -0.25, 1.875, 1.925, 2.5, 1.95, 1.85, -1, -1, 2.1, -1, 0.925, 0.95, -1)
$ws.Cells.Item(202, 5).Style = $ws.Cells.Item(200, 5).Style

$ws.Cells.Item(203, 1).Value = 201
$ws.Cells.Item(203, 1).Style = $ws.Cells.Item(200, 1).Style
Set-RowValues $ws "B203:AC203" @(7793515, "Chile Primera B", "Chile Primera B", 45397.79166666666, "Curico Unido", "Magallanes", 0, 0, "D", 2.5, 3.4, 2.375, 2.875, 3.5, 2.1, 0.25, 1.875, 1.925, 2.75, 1.8, 2, -1, 2.5, -1, 0.4375, -0.5, -1, 1)
$ws.Cells.Item(203, 5).Style = $ws.Cells.Item(200, 5).Style

